# Applies the coin-table refresh captured in the commit diff:
#  - price (column D) ticks for several existing rows
#  - a block of rows (9-17 and 41-43) shifted up by one rank, with the
#    row that fell off the top re-appearing at the bottom with a new price
#  - a couple of "Bestin24h" rank-change suffix moves in column E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B/C/E cells: plain text (coin name, link, rank+symbol string).
# These never look like numbers, so a normal .Value assignment keeps
# them stored as text, same as the source file.
$textUpdates = @{
    "B9" = "WazirX"
    "C9" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E9" = "8WazirXWRX"
    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E10" = "9MandalaExchangeTokenMDX"
    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E11" = "10LiechtensteinCryptoassetsExchangeLCX"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E12" = "11BitrueCoinBTR"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E13" = "12BitMartTokenBMX"
    "B14" = "MCDex"
    "C14" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E14" = "13MCDexMCB"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E15" = "14BitForexTokenBF"
    "B16" = "CoinExToken"
    "C16" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E16" = "15CoinExTokenCET"
    "B17" = "One"
    "C17" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E17" = "16OneONE"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E41" = "40BKEXTokenBKK"
    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "E42" = "41CEJICEJI"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "E43" = "42KickTokenKICK"
    "E47" = "46CoinbaseStockTokenCOINBestin24h"
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Column D cells: the "Price" column. These look numeric, and the source
# file stores them as literal text (e.g. "0.001600", "0.9610") where
# trailing/leading zeros must be preserved exactly. Force text number
# format first so Excel does not silently coerce them into doubles.
$priceUpdates = @{
    "D2" = "244.98"
    "D3" = "21.95"
    "D4" = "5.395"
    "D5" = "0.06005"
    "D7" = "0.8126"
    "D8" = "0.9610"
    "D9" = "0.1424"
    "D10" = "0.07388"
    "D11" = "0.03389"
    "D12" = "0.03055"
    "D13" = "0.09420"
    "D14" = "4.002"
    "D15" = "0.001600"
    "D16" = "0.04805"
    "D17" = "0.0005872"
    "D18" = "0.006217"
    "D19" = "0.005064"
    "D20" = "0.0009906"
    "D23" = "6.412"
    "D26" = "0.1292"
    "D40" = "0.04001"
    "D41" = "0.1074"
    "D42" = "0.002721"
    "D43" = "0.003019"
    "D44" = "0.005855"
    "D45" = "0.00005276"
    "D48" = "0.02177"
}

foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
}
